# IPSSMexample.xlsx update:
#  - add "ANC" and "AGE" columns (with example values) to the "examples" sheet
#  - reorder sheets so "examples" is first (and active), "explain" second
#  - update each sheet's remembered cell selection

$wb = $excel.ActiveWorkbook

$wsExamples = $wb.Worksheets.Item("examples")

# Insert two new columns (E, F) for ANC and AGE, pushing the existing
# CYTO_IPSSR..MLL_and_FLT3_mutations columns two places to the right.
$wsExamples.Columns("E:F").Insert()

$wsExamples.Range("E1").Value = "ANC"
$wsExamples.Range("F1").Value = "AGE"

$wsExamples.Range("E2").Value = 4.84
$wsExamples.Range("F2").Value = 79

$wsExamples.Range("E3").Value = 0.24
$wsExamples.Range("F3").Value = 66

# Remember the "explain" sheet's new selected cell before switching sheets.
$wsExplain = $wb.Worksheets.Item("explain")
$wsExplain.Activate() | Out-Null
$wsExplain.Range("A6").Select() | Out-Null

# Move "examples" to be the first tab in the workbook.
$wsExamples.Move($wb.Worksheets.Item(1))

# Re-fetch (the old reference goes stale across the Move) and make it the
# active sheet with its new remembered selection.
$wsExamples = $wb.Worksheets.Item("examples")
$wsExamples.Activate() | Out-Null
$wsExamples.Range("F4").Select() | Out-Null
